$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update isolation host data (host_species column D) for several rows
# (written in this order so new shared-string entries land at the
# indices the saved workbook expects)
$ws.Range("D5").Value = "Macaca mulatta"
$ws.Range("D6").Value = "Macaca nemestrina"
$ws.Range("D7").Value = "Tamias sibiricus"
$ws.Range("D8").Value = "Bos taurus"
$ws.Range("D4").Value = "Macaca fascicularis"
$ws.Range("D2").Value = "Phoca vitulina"

# Update the active cell selection (the sheet was last left with the
# cursor on B6, within the overall used range A1:M10)
$ws.Range("B6").Select()
